# Cronograma.xlsx update — "Checklist Qualidade e Cronograma finalizados"
#
# Strategy:
#  - Row 1 (header) is untouched (its text doesn't change, only an internal
#    shared-string index shift happens automatically).
#  - Rows 2-5 already have the right fill/font styling (alternating
#    light-green / white banding with bold-white header), we just correct a
#    few values (dates, a shortened "Ger. Requisitos" label, a renamed
#    ".xlsx" deliverable).
#  - Rows 6-14 are brand new data rows; we copy the formatting from the
#    existing banded rows (row 2 = "even" style, row 3 = "odd" style) with
#    Copy/PasteSpecial(formats) so the engine reuses the very same cellXfs
#    entries instead of inventing new ones, then fill in the values.
#  - Row 18 is a new single-line footer ("* Caso Haja Necessidade") with a
#    bold white-on-black cell (A18) followed by two bold white-on-white
#    filler cells (B18:C18) — built the same way: copy the bold header
#    format from A1, then recolor the fill via ThemeColor.
#  - Column B/C/F get widened/narrowed, and the view is rescrolled/rezoomed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix values in the existing rows 2-5 (styles already correct)
# ---------------------------------------------------------------------

# Row 2 — unchanged content, left as-is.

# Row 3 — Engenharia de Requisitos
$ws.Range("C3").Value2 = "Ger. Requisitos (Rogério A.), Analistas (Luiz E. , Thiago N.) "
$ws.Range("E3").Value2 = 43057

# Row 4 — Análise de Inconsistência
$ws.Range("D4").Value2 = 43057
$ws.Range("E4").Value2 = 43058
$ws.Range("F4").Value2 = "Checklist Inconsistência.xlsx"

# Row 5 — Rastreabilidade
$ws.Range("C5").Value2 = "Ger. Requisitos (Rogério A.), Analistas (Luiz E. , Thiago N.)"
$ws.Range("D5").Value2 = 43058
$ws.Range("E5").Value2 = 43058

# ---------------------------------------------------------------------
# 2. Build rows 6-14 (copy banding format from row 2 / row 3, then fill)
# ---------------------------------------------------------------------

function Set-BandedRow($rowNum, $templateRow, $a, $b, $c, $dDate, $eDate, $f) {
    $ws.Range("A$templateRow`:F$templateRow").Copy() | Out-Null
    $ws.Range("A$rowNum`:F$rowNum").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$rowNum").Value2 = $a
    $ws.Range("B$rowNum").Value2 = $b
    $ws.Range("C$rowNum").Value2 = $c
    $ws.Range("D$rowNum").Value2 = $dDate
    $ws.Range("E$rowNum").Value2 = $eDate
    if ($f -ne $null) {
        $ws.Range("F$rowNum").Value2 = $f
    }
}

Set-BandedRow 6  2 "Comprometimento Equipe" "Gerar Comprometimento da Equipe" "Gerente de Requisitos (Rogério A.)" 43059 43060 "Comprometimento da Equipe Técnica aos Requisitos.docx"
Set-BandedRow 7  3 "Validação com Forcedor" "Gerar Validação de Requisitos" "Gerente de Requisitos (Rogério A.)" 43060 43061 "Validação de Requisitos.docx"
Set-BandedRow 8  2 "Travar Baseline" "Gerar e Bloquear a Baseline" "Gerente do Projeto (Matheus F)" 43061 43061 "Baseline.docx"
Set-BandedRow 9  3 "Avaliação de Processo" "Avaliar Processo do Software" "Gerente de Qualidade (Marcus T)" 43061 43061 "Checklist Qualidade.xlsx"
Set-BandedRow 10 2 " Não Conformidade*" "Documentar Não Conformidade" "Gerente de Qualidade (Marcus T)" 43061 43061 "Não Conformidade.docx"
Set-BandedRow 11 3 "Ações Corretivas* " "Ações para Corrigir Não Confirmidades" "Gerente Responsável" 43061 43063 "Ações Corretivas.docx"
Set-BandedRow 12 2 "Acompanhar Tratamento* " "Andamento das Ações Corretivas" "Gerente de Qualidade (Marcus T)" 43061 43063 $null
Set-BandedRow 13 3 "Registrar Conclusão*" "Encerrar Não conformidades" "Gerente de Qualidade (Marcus T)" 43063 43063 "Conclusão de um Ação Corretiva.docx"
Set-BandedRow 14 2 "Finalizar Projeto" "Encerrar Projeto" "Gerente do Projeto (Matheus F)" 43063 43063 "GYN SOLUTION"

# ---------------------------------------------------------------------
# 3. Row 18 footer — bold white text, black fill (A18) / white fill (B18:C18)
# ---------------------------------------------------------------------

$ws.Range("A1").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Interior.ThemeColor = 1
$ws.Range("A18").Value2 = "* Caso Haja Necessidade"

$ws.Range("A1").Copy() | Out-Null
$ws.Range("B18:C18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18:C18").Interior.ThemeColor = 2

# ---------------------------------------------------------------------
# 4. Column widths
# ---------------------------------------------------------------------

$ws.Columns.Item(2).ColumnWidth = 35.6
$ws.Columns.Item(3).ColumnWidth = 51.42
$ws.Columns.Item(6).ColumnWidth = 51.26

# ---------------------------------------------------------------------
# 5. Sheet view — selection + zoom
# ---------------------------------------------------------------------

$ws.Range("B23").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
